$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 14212.5
$ws.Range("J17").Value = 16000
$ws.Range("L17").Value = 48000
$ws.Range("N17").Value = -48336

$ws.Range("H112").Value = 1920.3077
$ws.Range("J112").Value = 1920.3077
$ws.Range("L112").Value = 5760.9231
$ws.Range("N112").Value = -7976.9231

$ws.Range("H128").Value = 60776
$ws.Range("J128").Value = 60776
$ws.Range("L128").Value = 60776
$ws.Range("N128").Value = -70736

$ws.Range("H138").Value = 3617.551
$ws.Range("I138").Value = 1494.5
$ws.Range("J138").Value = 4850.2905
$ws.Range("K138").Value = 4483.5
$ws.Range("L138").Value = 14550.8715
$ws.Range("M138").Value = 656.5
$ws.Range("N138").Value = -24830.8715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5411.5
$ws.Range("I2").Value = 379.44446
$ws.Range("K2").Value = 379.44446
$ws.Range("M2").Value = -266.44446

$ws.Range("H16").Value = 690.1429000000001
$ws.Range("I16").Value = 526.2
$ws.Range("K16").Value = 526.2
$ws.Range("M16").Value = -239.2

$ws.Range("H32").Value = 5002.101
$ws.Range("I32").Value = 4093.117
$ws.Range("K32").Value = 4093.117
$ws.Range("M32").Value = -3806.117

$ws.Range("H35").Value = 1500
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 1500
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 1500
$ws.Range("N35").Value = -2312
$ws.Range("M35").ClearContents()

$ws.Range("H116").Value = 5411.5
$ws.Range("I116").Value = 379.44446
$ws.Range("K116").Value = 379.44446
$ws.Range("M116").Value = 1914.55554

$ws.Range("H122").Value = 2028.9836
$ws.Range("I122").Value = 1634.262
$ws.Range("K122").Value = 4902.786
$ws.Range("M122").Value = -2452.786

$ws.Range("H132").Value = 2421.5938
$ws.Range("I132").Value = 967.08
$ws.Range("K132").Value = 2901.24
$ws.Range("M132").Value = -371.2400000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5411.5
$ws.Range("I3").Value = 379.44446
$ws.Range("K3").Value = 379.44446
$ws.Range("M3").Value = -265.44446

$ws.Range("H25").Value = 18000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 18000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 18000
$ws.Range("N25").Value = -18470
$ws.Range("M25").ClearContents()

$ws.Range("H86").Value = 5861.316
$ws.Range("I86").Value = 5400.9062
$ws.Range("K86").Value = 5400.9062
$ws.Range("M86").Value = -4277.9062

$ws.Range("H89").Value = 5861.316
$ws.Range("I89").Value = 5400.9062
$ws.Range("K89").Value = 27004.531
$ws.Range("M89").Value = -21388.531

$ws.Range("H105").Value = 15651.792
$ws.Range("I105").Value = 12603.474
$ws.Range("K105").Value = 12603.474
$ws.Range("M105").Value = -10856.474

$ws.Range("H134").Value = 1660.1702
$ws.Range("I134").Value = 1109.5122
$ws.Range("J134").Value = 5423
$ws.Range("K134").Value = 3328.536599999999
$ws.Range("L134").Value = 16269
$ws.Range("M134").Value = -793.5365999999995
$ws.Range("N134").Value = -21339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 9035.727999999999
$ws.Range("J22").Value = 10149.167
$ws.Range("L22").Value = 10149.167
$ws.Range("N22").Value = -10849.167

$ws.Range("H33").Value = 1500
$ws.Range("I33").Value = 1500
$ws.Range("K33").Value = 1500
$ws.Range("M33").Value = -1121

$ws.Range("H51").Value = 24995
$ws.Range("I51").Value = 24995
$ws.Range("K51").Value = 24995
$ws.Range("M51").Value = -24259

$ws.Range("H61").Value = 24995
$ws.Range("I61").Value = 24995
$ws.Range("K61").Value = 24995
$ws.Range("M61").Value = -24647

$ws.Range("H68").Value = 69332
$ws.Range("J68").Value = 69332
$ws.Range("L68").Value = 69332
$ws.Range("N68").Value = -70830

$ws.Range("H71").Value = 69332
$ws.Range("J71").Value = 69332
$ws.Range("L71").Value = 207996
$ws.Range("N71").Value = -215484

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 844.625
$ws.Range("J5").Value = 5562222
$ws.Range("K5").Value = 2533.875
$ws.Range("L5").Value = 16686666
$ws.Range("M5").Value = -2421.875
$ws.Range("N5").Value = -16686890

$ws.Range("H12").Value = 176.33333
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 176.33333
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 528.99999
$ws.Range("N12").Value = -874.99999
$ws.Range("M12").ClearContents()

$ws.Range("H16").Value = 258
$ws.Range("I16").Value = 137.5
$ws.Range("J16").Value = 499
$ws.Range("K16").Value = 412.5
$ws.Range("L16").Value = 1497
$ws.Range("M16").Value = -239.5
$ws.Range("N16").Value = -1843

$ws.Range("H34").Value = 8646664
$ws.Range("J34").Value = 6999.75
$ws.Range("L34").Value = 20999.25
$ws.Range("N34").Value = -21167.25

$ws.Range("H55").Value = 1745.5
$ws.Range("I55").Value = 1495.1538
$ws.Range("K55").Value = 4485.4614
$ws.Range("M55").Value = -4308.4614

$ws.Range("H131").Value = 6265272
$ws.Range("I131").Value = 8929640
$ws.Range("K131").Value = 26788920
$ws.Range("M131").Value = -26783880

$ws.Range("I135").Value = 844.625
$ws.Range("J135").Value = 5562222
$ws.Range("K135").Value = 7601.625
$ws.Range("L135").Value = 50059998
$ws.Range("M135").Value = -5066.625
$ws.Range("N135").Value = -50065068

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 3151870.5
$ws.Range("I12").Value = 4398
$ws.Range("J12").Value = 5512475
$ws.Range("K12").Value = 4398
$ws.Range("L12").Value = 5512475
$ws.Range("M12").Value = -4258
$ws.Range("N12").Value = -5512755

$ws.Range("H113").Value = 2538.75
$ws.Range("I113").Value = 1858.6957
$ws.Range("K113").Value = 1858.6957
$ws.Range("M113").Value = 311.3043

$ws.Range("H126").Value = 2330.762
$ws.Range("I126").Value = 1411.4828
$ws.Range("K126").Value = 4234.4484
$ws.Range("M126").Value = -1764.4484

$ws.Range("H128").Value = 55000
$ws.Range("J128").Value = 55000
$ws.Range("L128").Value = 55000
$ws.Range("N128").Value = -64960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6105.885
$ws.Range("I7").Value = 3867.5217
$ws.Range("K7").Value = 3867.5217
$ws.Range("M7").Value = -3755.5217

$ws.Range("H20").Value = 4482
$ws.Range("J20").Value = 4482
$ws.Range("L20").Value = 4482
$ws.Range("N20").Value = -4934

$ws.Range("H22").Value = 3745.95
$ws.Range("I22").Value = 867.6667
$ws.Range("K22").Value = 867.6667
$ws.Range("M22").Value = -572.6667

$ws.Range("H27").Value = 3745.95
$ws.Range("I27").Value = 867.6667
$ws.Range("K27").Value = 867.6667
$ws.Range("M27").Value = -760.6667

$ws.Range("H31").Value = 8403.817999999999
$ws.Range("I31").Value = 276.7143
$ws.Range("J31").Value = 22626.25
$ws.Range("K31").Value = 276.7143
$ws.Range("L31").Value = 22626.25
$ws.Range("M31").Value = -28.71429999999998
$ws.Range("N31").Value = -23122.25

$ws.Range("H126").Value = 6105.885
$ws.Range("I126").Value = 3867.5217
$ws.Range("K126").Value = 11602.5651
$ws.Range("M126").Value = -9132.5651

$ws.Range("H132").Value = 3831.5
$ws.Range("I132").Value = 3318.7932
$ws.Range("K132").Value = 9956.3796
$ws.Range("M132").Value = -7426.3796

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 30000
$ws.Range("J33").Value = 30000
$ws.Range("L33").Value = 30000
$ws.Range("N33").Value = -30500

$ws.Range("H36").Value = 30000
$ws.Range("J36").Value = 30000
$ws.Range("L36").Value = 30000
$ws.Range("N36").Value = -30500

$ws.Range("H122").Value = 2457.7646
$ws.Range("I122").Value = 1780.4615
$ws.Range("K122").Value = 5341.3845
$ws.Range("M122").Value = -2891.3845

$ws.Range("H132").Value = 3175.6924
$ws.Range("I132").Value = 2542.52
$ws.Range("K132").Value = 7627.559999999999
$ws.Range("M132").Value = -5097.559999999999
